$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5195787887265466
$ws.Range("C3").Value = 0.1899731131478304
$ws.Range("C8").Value = 0.308413692671376
$ws.Range("C9").Value = 0.1734898436308585
$ws.Range("C10").Value = 0.2055719433594982
